$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Status changed from "Ready for handoff" to "Handed back: in sync with en-US"
# everywhere it appears: the Overview summary columns for each language, and
# each language sheet's own Status column.
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("C2").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("C2").Value = "Handed back: in sync with en-US"

# zh-cn sheet: refresh handback datetime and clear the stale-handback error detail
$wsZhCn.Range("K2").Value = "2016-08-15 14:47:27"
$wsZhCn.Range("P2").Value = ""

# de-de sheet: refresh handback datetime and clear the stale-handback error detail
$wsDeDe.Range("K2").Value = "2016-08-15 14:47:34"
$wsDeDe.Range("P2").Value = ""
